# Reference.xlsx edit
# Commit: Shortened paths. Optimize/Macro/../Library -> Optimize/Library.
#         Cases where we have one descriptor for each NRG, went to 1 descriptor for All_NRG
#
# This consolidates the "one growth-rate descriptor per NRG" row (old row 6,
# "Coal_Grow_2_1-Gas_Grow_2_1-..." / the stray "`" placeholder) into a single
# "Debug" analysis case that uses one shared "All_NRG_Grow_2_1" descriptor,
# and adds a new "5% Cheaper Per Year" case (row 7) plus a second FileName
# entry reusing Nuke_Cap$_0.5_1 for the new Debug/Half_Cap combo (row 8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: new "Debug" Analysis case (was blank Analysis cols, stray "`") ---
$ws.Range("M6").Value = "Debug"
$ws.Range("N6").Value = "Debug"
$ws.Range("O6").Formula = "=M6"
$ws.Range("P6").Value = "Double_All"
$ws.Range("Q6").Value = "By 2x"

# --- Row 7: new "5%_Per_Year" Case entry ---
$ws.Range("P7").Value = "5%_Per_Year"
$ws.Range("Q7").Value = "5% Cheaper Per Year"
$ws.Range("R7").Formula = "=P7"
$ws.Range("S7").Formula = '=M4 & "_" &P7'
$ws.Range("T7").Value = 'Nuke_Cap$_1_0.95'

# --- Row 6 FileName: now a single "All_NRG" descriptor instead of one per NRG ---
$ws.Range("T6").Value = "All_NRG_Grow_2_1"

# --- Row 8: new Analysis_Case / FileName entry for Debug + Half_Cap ---
$ws.Range("S8").Formula = '=M6 & "_" & P5'
$ws.Range("T8").Value = 'Nuke_Cap$_0.5_1'

# --- Row heights follow the new wrapped-text content ---
$ws.Rows.Item(6).RowHeight = 29.15
$ws.Rows.Item(7).RowHeight = 43.9
$ws.Rows.Item(8).RowHeight = 29.15
$ws.Rows.Item(9).RowHeight = 14.65

# --- Selection moves to U7 ---
$ws.Range("U7").Select()

$excel.Calculate()
